$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "Paso 2. Análisis del Código" - second bullet ("Debe desglosar...")
#    is emptied out: keep the "Prrafodelista" paragraph style but drop
#    the list numbering and the run text entirely.
# ---------------------------------------------------------------------
$descBullet = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs($i)
    if ($cand.Range.Text -like "Debe desglosar*") {
        $descBullet = $cand
        break
    }
}
if ($descBullet -ne $null) {
    $descBullet.Range.Delete()
}

# ---------------------------------------------------------------------
# 2) "Paso 3. Experimentación con Parámetros" - first bullet is reworded
#    from:
#      "Debe modificar los parámetros clave del controlador
#       proporcional, como la constante proporcional (Kp)."
#    to:
#      "Debe modificar la constante proporcional Kp con un rango de 2 a
#       10, utilice serialplot para identificar qué ocurre con el
#       sistema."
#    The "Kp" run (wrapped in proofErr spell-check tags) is left alone
#    so only the text immediately before and after it is touched.
# ---------------------------------------------------------------------
$modBullet = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs($i)
    if ($cand.Range.Text -like "Debe modificar los parámetros*") {
        $modBullet = $cand
        break
    }
}
if ($modBullet -ne $null) {
    $rngA = $modBullet.Range
    $rngA.Find.Execute(
        "Debe modificar los parámetros clave del controlador proporcional, como la constante proporcional (",
        $true, $false, $false, $false, $false, $true, 1, $false,
        "Debe modificar la constante proporcional ", 2) | Out-Null

    $modBullet2 = $null
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $cand = $d.Paragraphs($i)
        if ($cand.Range.Text -like "Debe modificar la constante proporcional *") {
            $modBullet2 = $cand
            break
        }
    }
    $rngB = $modBullet2.Range
    $rngB.Find.Execute(
        ").", $true, $false, $false, $false, $false, $true, 1, $false,
        " con un rango de 2 a 10, utilice serialplot para identificar qué ocurre con el sistema.", 2) | Out-Null
}

# ---------------------------------------------------------------------
# 3) "Recolección de Datos" bullet - reword "Debe documentar..." and
#    promote it into the bulleted ("Prrafodelista") list, then add two
#    new bullets after it:
#      - "Debe responder a las siguientes preguntas: ..."
#      - an empty "Prrafodelista" spacer paragraph.
# ---------------------------------------------------------------------
$docBullet = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs($i)
    if ($cand.Range.Text -like "Debe documentar los resultados*") {
        $docBullet = $cand
        break
    }
}

# Paragraph that already carries the plain "-" bulleted list style we
# want to reuse (the "Objetivos Específicos" bullets).
$listSample = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs($i)
    if ($cand.Range.Text -like "Analizar el código proporcionado del controlador*") {
        $listSample = $cand
        break
    }
}

if ($docBullet -ne $null) {
    $rngC = $docBullet.Range
    $rngC.Find.Execute(
        "Debe documentar los resultados de las pruebas, registrando las configuraciones específicas de los parámetros y observando el comportamiento del carrito en términos de velocidad, estabilidad y capacidad para evadir obstáculos.",
        $true, $false, $false, $false, $false, $true, 1, $false,
        "Debe documentar los resultados de las pruebas, registrando el comportamiento específico de los valores seleccionados Kp dentro de los rangos solicitados.", 2) | Out-Null

    $docBullet2 = $null
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $cand = $d.Paragraphs($i)
        if ($cand.Range.Text -like "Debe documentar los resultados de las pruebas, registrando el comportamiento*") {
            $docBullet2 = $cand
            break
        }
    }

    # Add the new question bullet right after it first.
    $docBullet2.Range.InsertParagraphAfter()
    $questionPara = $docBullet2.Next()
    $questionPara.Range.InsertBefore("Debe responder a las siguientes preguntas: ¿Qué ocurre con el robot cuando el kp es demasiado alto? ¿Y cuándo es muy bajo? Justifique su respuesta.")

    # Apply the bullet list style to the new question paragraph first...
    $questionPara2 = $docBullet2.Next()
    $questionPara2.Style = "Prrafodelista"
    $questionPara2.Range.ListFormat.ApplyListTemplate($listSample.Range.ListFormat.ListTemplate) | Out-Null

    # ...then to the reworded "Debe documentar" paragraph.
    $docBullet3 = $questionPara2.Previous()
    $docBullet3.Style = "Prrafodelista"
    $docBullet3.Range.ListFormat.ApplyListTemplate($listSample.Range.ListFormat.ListTemplate) | Out-Null

    # Empty spacer paragraph after the question bullet, still styled as
    # "Prrafodelista" but without any numbering.
    $questionPara3 = $docBullet3.Next()
    $questionPara3.Range.InsertParagraphAfter()
    $spacerPara = $questionPara3.Next()
    $spacerPara.Range.Style = "Prrafodelista"
}

Write-Output "done"
